# Weekly update to "Fruta, Agrícola del Norte S.A. de Arica - Damasco":
# the price records in rows 2-13 are rotated to reflect the latest
# weekly refresh (row 6 keeps its original data; all others are
# replaced by another row's previous values per the source diff).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 now holds the data previously in row 13
$ws.Cells.Item(2, 4).Value = 44894
$ws.Cells.Item(2, 12).Value = 'Segunda'
$ws.Cells.Item(2, 13).Value = 130
$ws.Cells.Item(2, 14).Value = 19000
$ws.Cells.Item(2, 15).Value = 20000
$ws.Cells.Item(2, 16).Value = 19462
$ws.Cells.Item(2, 17).Value = '$/caja 16 kilos granel'
$ws.Cells.Item(2, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(2, 19).Value = 1216
$ws.Cells.Item(2, 20).Value = 16

# Row 3 now holds the data previously in row 7
$ws.Cells.Item(3, 4).Value = 44533
$ws.Cells.Item(3, 12).Value = 'Primera'
$ws.Cells.Item(3, 13).Value = 140
$ws.Cells.Item(3, 14).Value = 14000
$ws.Cells.Item(3, 15).Value = 15000
$ws.Cells.Item(3, 16).Value = 14500
$ws.Cells.Item(3, 17).Value = '$/caja 10 kilos'
$ws.Cells.Item(3, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(3, 19).Value = 1450
$ws.Cells.Item(3, 20).Value = 10

# Row 4 now holds the data previously in row 8
$ws.Cells.Item(4, 4).Value = 44901
$ws.Cells.Item(4, 12).Value = 'Segunda'
$ws.Cells.Item(4, 13).Value = 200
$ws.Cells.Item(4, 14).Value = 17000
$ws.Cells.Item(4, 15).Value = 18000
$ws.Cells.Item(4, 16).Value = 17500
$ws.Cells.Item(4, 17).Value = '$/bandeja 18 kilos'
$ws.Cells.Item(4, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(4, 19).Value = 972
$ws.Cells.Item(4, 20).Value = 18

# Row 5 now holds the data previously in row 4
$ws.Cells.Item(5, 4).Value = 44881
$ws.Cells.Item(5, 12).Value = 'Segunda'
$ws.Cells.Item(5, 13).Value = 300
$ws.Cells.Item(5, 14).Value = 41000
$ws.Cells.Item(5, 15).Value = 42000
$ws.Cells.Item(5, 16).Value = 41500
$ws.Cells.Item(5, 17).Value = '$/bandeja 18 kilos'
$ws.Cells.Item(5, 18).Value = 'Región de Coquimbo'
$ws.Cells.Item(5, 19).Value = 2306
$ws.Cells.Item(5, 20).Value = 18

# Row 6: unchanged
# Row 7 now holds the data previously in row 3
$ws.Cells.Item(7, 4).Value = 44545
$ws.Cells.Item(7, 12).Value = 'Primera'
$ws.Cells.Item(7, 13).Value = 200
$ws.Cells.Item(7, 14).Value = 24000
$ws.Cells.Item(7, 15).Value = 25000
$ws.Cells.Item(7, 16).Value = 24500
$ws.Cells.Item(7, 17).Value = '$/bandeja 18 kilos'
$ws.Cells.Item(7, 18).Value = 'Región de Coquimbo'
$ws.Cells.Item(7, 19).Value = 1361
$ws.Cells.Item(7, 20).Value = 18

# Row 8 now holds the data previously in row 11
$ws.Cells.Item(8, 4).Value = 44895
$ws.Cells.Item(8, 12).Value = 'Segunda'
$ws.Cells.Item(8, 13).Value = 130
$ws.Cells.Item(8, 14).Value = 19000
$ws.Cells.Item(8, 15).Value = 20000
$ws.Cells.Item(8, 16).Value = 19462
$ws.Cells.Item(8, 17).Value = '$/caja 16 kilos granel'
$ws.Cells.Item(8, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(8, 19).Value = 1216
$ws.Cells.Item(8, 20).Value = 16

# Row 9 now holds the data previously in row 5
$ws.Cells.Item(9, 4).Value = 44524
$ws.Cells.Item(9, 12).Value = 'Segunda'
$ws.Cells.Item(9, 13).Value = 200
$ws.Cells.Item(9, 14).Value = 27000
$ws.Cells.Item(9, 15).Value = 28000
$ws.Cells.Item(9, 16).Value = 27500
$ws.Cells.Item(9, 17).Value = '$/bandeja 18 kilos'
$ws.Cells.Item(9, 18).Value = 'Provincia de San Felipe de Aconcagua'
$ws.Cells.Item(9, 19).Value = 1528
$ws.Cells.Item(9, 20).Value = 18

# Row 10 now holds the data previously in row 2
$ws.Cells.Item(10, 4).Value = 44917
$ws.Cells.Item(10, 12).Value = 'Segunda'
$ws.Cells.Item(10, 13).Value = 250
$ws.Cells.Item(10, 14).Value = 20000
$ws.Cells.Item(10, 15).Value = 23000
$ws.Cells.Item(10, 16).Value = 21800
$ws.Cells.Item(10, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(10, 18).Value = 'Región de Coquimbo'
$ws.Cells.Item(10, 19).Value = 1211
$ws.Cells.Item(10, 20).Value = 18

# Row 11 now holds the data previously in row 12
$ws.Cells.Item(11, 4).Value = 44160
$ws.Cells.Item(11, 12).Value = 'Primera'
$ws.Cells.Item(11, 13).Value = 250
$ws.Cells.Item(11, 14).Value = 24000
$ws.Cells.Item(11, 15).Value = 25000
$ws.Cells.Item(11, 16).Value = 24500
$ws.Cells.Item(11, 17).Value = '$/bandeja 18 kilos'
$ws.Cells.Item(11, 18).Value = 'Provincia de San Felipe de Aconcagua'
$ws.Cells.Item(11, 19).Value = 1361
$ws.Cells.Item(11, 20).Value = 18

# Row 12 now holds the data previously in row 9
$ws.Cells.Item(12, 4).Value = 44880
$ws.Cells.Item(12, 12).Value = 'Primera'
$ws.Cells.Item(12, 13).Value = 200
$ws.Cells.Item(12, 14).Value = 33000
$ws.Cells.Item(12, 15).Value = 34000
$ws.Cells.Item(12, 16).Value = 33500
$ws.Cells.Item(12, 17).Value = '$/caja 10 kilos'
$ws.Cells.Item(12, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(12, 19).Value = 3350
$ws.Cells.Item(12, 20).Value = 10

# Row 13 now holds the data previously in row 10
$ws.Cells.Item(13, 4).Value = 44174
$ws.Cells.Item(13, 12).Value = 'Primera'
$ws.Cells.Item(13, 13).Value = 300
$ws.Cells.Item(13, 14).Value = 19000
$ws.Cells.Item(13, 15).Value = 20000
$ws.Cells.Item(13, 16).Value = 19500
$ws.Cells.Item(13, 17).Value = '$/bandeja 18 kilos'
$ws.Cells.Item(13, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(13, 19).Value = 1083
$ws.Cells.Item(13, 20).Value = 18

